$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows before the existing row 166 (old rows 166-190 shift down to 171-195).
$ws.Range("A166:A170").EntireRow.Insert()

# Row 166 - Mandarina, Murcott, Primera (2021-10-05)
$ws.Cells.Item(166, 1).Value = 5
$ws.Cells.Item(166, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(166, 3).Value = 'Maule'
$ws.Cells.Item(166, 4).Value = '2021-10-05'
$ws.Cells.Item(166, 5).Value = 7
$ws.Cells.Item(166, 6).Value = 'Fruta'
$ws.Cells.Item(166, 7).Value = 100102
$ws.Cells.Item(166, 8).Value = 'Cítricos'
$ws.Cells.Item(166, 9).Value = 100102004
$ws.Cells.Item(166, 10).Value = 'Mandarina'
$ws.Cells.Item(166, 11).Value = 'Murcott'
$ws.Cells.Item(166, 12).Value = 'Primera'
$ws.Cells.Item(166, 13).Value = 300
$ws.Cells.Item(166, 14).Value = 7000
$ws.Cells.Item(166, 15).Value = 7000
$ws.Cells.Item(166, 16).Value = 7000
$ws.Cells.Item(166, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(166, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(166, 19).Value = 389
$ws.Cells.Item(166, 20).Value = 18

# Row 167 - Mandarina, Murcott, Primera (2021-10-05)
$ws.Cells.Item(167, 1).Value = 5
$ws.Cells.Item(167, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(167, 3).Value = 'Maule'
$ws.Cells.Item(167, 4).Value = '2021-10-05'
$ws.Cells.Item(167, 5).Value = 7
$ws.Cells.Item(167, 6).Value = 'Fruta'
$ws.Cells.Item(167, 7).Value = 100102
$ws.Cells.Item(167, 8).Value = 'Cítricos'
$ws.Cells.Item(167, 9).Value = 100102004
$ws.Cells.Item(167, 10).Value = 'Mandarina'
$ws.Cells.Item(167, 11).Value = 'Murcott'
$ws.Cells.Item(167, 12).Value = 'Primera'
$ws.Cells.Item(167, 13).Value = 350
$ws.Cells.Item(167, 14).Value = 6000
$ws.Cells.Item(167, 15).Value = 6000
$ws.Cells.Item(167, 16).Value = 6000
$ws.Cells.Item(167, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(167, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(167, 19).Value = 333
$ws.Cells.Item(167, 20).Value = 18

# Row 168 - Mandarina, Murcott, Primera (2021-10-05)
$ws.Cells.Item(168, 1).Value = 5
$ws.Cells.Item(168, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(168, 3).Value = 'Maule'
$ws.Cells.Item(168, 4).Value = '2021-10-05'
$ws.Cells.Item(168, 5).Value = 7
$ws.Cells.Item(168, 6).Value = 'Fruta'
$ws.Cells.Item(168, 7).Value = 100102
$ws.Cells.Item(168, 8).Value = 'Cítricos'
$ws.Cells.Item(168, 9).Value = 100102004
$ws.Cells.Item(168, 10).Value = 'Mandarina'
$ws.Cells.Item(168, 11).Value = 'Murcott'
$ws.Cells.Item(168, 12).Value = 'Primera'
$ws.Cells.Item(168, 13).Value = 350
$ws.Cells.Item(168, 14).Value = 6000
$ws.Cells.Item(168, 15).Value = 6000
$ws.Cells.Item(168, 16).Value = 6000
$ws.Cells.Item(168, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(168, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(168, 19).Value = 333
$ws.Cells.Item(168, 20).Value = 18

# Row 169 - Mandarina, Murcott, Segunda (2021-10-05)
$ws.Cells.Item(169, 1).Value = 5
$ws.Cells.Item(169, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(169, 3).Value = 'Maule'
$ws.Cells.Item(169, 4).Value = '2021-10-05'
$ws.Cells.Item(169, 5).Value = 7
$ws.Cells.Item(169, 6).Value = 'Fruta'
$ws.Cells.Item(169, 7).Value = 100102
$ws.Cells.Item(169, 8).Value = 'Cítricos'
$ws.Cells.Item(169, 9).Value = 100102004
$ws.Cells.Item(169, 10).Value = 'Mandarina'
$ws.Cells.Item(169, 11).Value = 'Murcott'
$ws.Cells.Item(169, 12).Value = 'Segunda'
$ws.Cells.Item(169, 13).Value = 200
$ws.Cells.Item(169, 14).Value = 3500
$ws.Cells.Item(169, 15).Value = 3500
$ws.Cells.Item(169, 16).Value = 3500
$ws.Cells.Item(169, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(169, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(169, 19).Value = 194
$ws.Cells.Item(169, 20).Value = 18

# Row 170 - Mandarina, Murcott, Tercera (2021-10-05)
$ws.Cells.Item(170, 1).Value = 5
$ws.Cells.Item(170, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(170, 3).Value = 'Maule'
$ws.Cells.Item(170, 4).Value = '2021-10-05'
$ws.Cells.Item(170, 5).Value = 7
$ws.Cells.Item(170, 6).Value = 'Fruta'
$ws.Cells.Item(170, 7).Value = 100102
$ws.Cells.Item(170, 8).Value = 'Cítricos'
$ws.Cells.Item(170, 9).Value = 100102004
$ws.Cells.Item(170, 10).Value = 'Mandarina'
$ws.Cells.Item(170, 11).Value = 'Murcott'
$ws.Cells.Item(170, 12).Value = 'Tercera'
$ws.Cells.Item(170, 13).Value = 240
$ws.Cells.Item(170, 14).Value = 4000
$ws.Cells.Item(170, 15).Value = 4000
$ws.Cells.Item(170, 16).Value = 4000
$ws.Cells.Item(170, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(170, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(170, 19).Value = 222
$ws.Cells.Item(170, 20).Value = 18
